# Align WHO Hepatitis B "3-dose series" sheet with current WHO/CDC guidance:
#  - Dose 3 absolute minimum age: 14 weeks -> 24 weeks
#  - Dose 3 preferable interval from Dose 2: 4 weeks -> 8 weeks
#  - Add a new "Preferable Interval" row: 16 weeks minimum from Dose 1 to Dose 3
#  - Drop now-unused trailing "n/a" filler cells left over from the old layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-dose series")

# --- Dose 3 / Age (row 27): absolute minimum age 14 weeks -> 24 weeks ---
$ws.Range("B27").Value2 = "24 weeks"
$ws.Range("C27").Value2 = "24 weeks"
$ws.Range("D27").Value2 = "24 weeks"

# Drop the unused trailing "n/a" columns on rows 8, 17 and 27 (Age rows)
$ws.Range("G8:H8").Clear()
$ws.Range("G17:H17").Clear()
$ws.Range("G27:H27").Clear()

# Drop the unused trailing "n/a" columns on rows 18 and 28 (Preferable Interval rows)
$ws.Range("J18:L18").Clear()
$ws.Range("J28:L28").Clear()

# --- Dose 3 / Preferable Interval from Dose 2 (row 28): 4 weeks -> 8 weeks ---
$ws.Range("F28").Value2 = "8 weeks"
$ws.Range("G28").Value2 = "8 weeks"

# --- Insert a new Preferable Interval row: 16-week minimum from Dose 1 to Dose 3 ---
$ws.Rows("29:29").Insert()
$ws.Range("A29").Value2 = "Preferable Interval"
$ws.Range("B29").Value2 = "n/a"
$ws.Range("C29").Value2 = "Dose 1"
$ws.Range("D29").Value2 = "n/a"
$ws.Range("E29").Value2 = "n/a"
$ws.Range("F29").Value2 = "16 weeks"
$ws.Range("G29").Value2 = "16 weeks"
$ws.Range("H29").Value2 = "16 weeks"
$ws.Range("I29").Value2 = "n/a"
$ws.Range("J29:L29").Clear()
